# Update "想去人数" (column F) figures on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 1079
    5  = 84
    7  = 55
    8  = 11186
    9  = 4285
    10 = 25
    11 = 22
    13 = 2495
    14 = 1068
    15 = 98
    17 = 157
    18 = 484
    19 = 11227
    20 = 11072
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
